# Add a new bullet point under "Principal Cloud and Platform Engineer
# (June 2020—January 2024)" describing the COVID-19 remote-learning work,
# matching the existing "Compact" bulleted-list style (numId 1002) used
# by its sibling bullets.

$d = $word.ActiveDocument

$anchorText = "Partnered with Enterprise Architecture"
$newBulletText = "As every school in America transitioned to online learning during the COVID-19 lockdowns, I was the technical/development lead on the team who supported all SRE and product engineering teams, working on core platforms and services."

# Locate the first bullet of the "Principal Cloud and Platform Engineer"
# job entry; the new bullet will be inserted immediately before it so it
# becomes the new first bullet of that job.
$findRange = $d.Content
$found = $findRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph '$anchorText'"
}
$findRange.Collapse(1)
$anchorPara = $findRange.Paragraphs(1)
$insertionPoint = $anchorPara.Range.Start

# Insert a new (currently empty) paragraph before the anchor paragraph.
# Because the anchor paragraph already carries the "Compact" style and
# numId=1002 list numbering, the freshly split-off paragraph inherits the
# exact same pPr (style + numPr), matching the sibling bullets exactly.
$anchorPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Range($insertionPoint, $insertionPoint).Paragraphs(1)
$newParaRange = $newPara.Range

# Replace the new (empty) paragraph's contents via raw OOXML so the
# resulting run carries xml:space="preserve" like every other run in
# this document.
$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1002"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $newBulletText + '</w:t></w:r></w:p>'
$newParaRange.InsertXML($xmlFragment) | Out-Null

Write-Host "Inserted new bullet under 'Principal Cloud and Platform Engineer'."
